$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 8 (pushes the "extr*" contingency rows down by 2),
# making room for two new line entries: line7 and line8.
$ws.Rows("8:9").Insert()

# New row 8 -> line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New row 9 -> line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Match formatting of the style used by the other "A" column index cells
# (bold font, thin box border, centered alignment).
$ws.Cells.Item(7, 1).Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# The contingency results for extr4 and extr5 (now at rows 13 and 14 after the
# insert) flip to "in service" = TRUE.
$ws.Cells.Item(13, 5).Value = $true
$ws.Cells.Item(14, 5).Value = $true

# Column A is a simple 0-based row index; renumber the shifted rows so it
# stays sequential after the insert.
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
